# testprep: create multiple tables for quiz and modeltest in the same DB
# Rename the five raw-data sheets from the generic "SheetN" defaults to
# "quizN" so each table is clearly identified.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet1").Name = "quiz1"
$wb.Worksheets.Item("Sheet2").Name = "quiz2"
$wb.Worksheets.Item("Sheet3").Name = "quiz3"
$wb.Worksheets.Item("Sheet4").Name = "quiz4"
$wb.Worksheets.Item("Sheet5").Name = "quiz5"

# Touch quiz2 (cell B1) the way Excel itself records the last selection on a
# sheet, then return focus to quiz1 so it stays the active tab.
$ws2 = $wb.Worksheets.Item("quiz2")
$ws2.Activate()
$ws2.Range("B1").Select()

$ws1 = $wb.Worksheets.Item("quiz1")
$ws1.Activate()
